$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: relocate the current row-4 question (the VCB / alpha/beta one) down to row 6 ---
# Row 6 currently only has A6:C6 and J6 filled in; D6:I6 are empty and need both the
# values AND the formatting (D4 has wrapText style) that used to live on row 4.
$ws.Range("D4:I4").Copy()
$ws.Range("D6").PasteSpecial(-4122)   # xlPasteFormats - bring the wrap-text style (and any other fmt) across first
$ws.Range("D4:I4").Copy()
$ws.Range("D6").PasteSpecial(-4163)   # xlPasteValues - now bring the actual text/values across
$excel.CutCopyMode = $false

# Row 6 becomes the tall row (used to be row 4's height); row 4 goes back to the default height.
$ws.Rows(6).RowHeight = 193.85
$ws.Rows(4).AutoFit()

# --- Step 2: write the brand-new question into row 4 ---
$ws.Range("D4").Value = 'Tính \(\lim\limits_{x \to 0} \frac{\sin 2x + \arcsin^2 x - \arctan^2 x}{3x}\)'
$ws.Range("E4").Value = ' +\infty'
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = ' \frac{2}{3} '
$ws.Range("H4").Value = ' \frac{1}{2}'
$ws.Range("I4").Value = ' \frac{2}{3} '
$ws.Range("D4").Style = "Normal"

$ws.Range("G4").NumberFormat = "d-mmm"
$ws.Range("I4").NumberFormat = "d-mmm"

# --- Step 3: row 5 only gains the "Correct Answer" entry in I5 ---
$ws.Range("I5").Value = 0

# --- Step 4: fill in the brand-new questions for rows 7-13 ---
$ws.Range("D7").Value = 'Tính \(\lim_{x \to 3} \frac{\sqrt{x + 1} - 2}{x - 3}\)'
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = ' \frac{1}{4}'
$ws.Range("G7").Value = '\infty'
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = ' \frac{1}{4}'

$ws.Range("D8").Value = 'Vô cùng bé nào sau đây tương đương với: \(\alpha(x) = \sin x - \tan x + x^3\)'
$ws.Range("E8").Value = 'x^3'
$ws.Range("F8").Value = ' \frac{x^3}{2} '
$ws.Range("G8").Value = ' \frac{x^2}{2} '
$ws.Range("H8").Value = 'x'
$ws.Range("I8").Value = ' \frac{x^3}{2} '

$ws.Range("D9").Value = 'Kết quả của giới hạn \(\lim\limits_{x \to 0} \frac{\sin 2x + 3 \cos x + x}{2x + \cos^2 3x}\) là:'
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 3

$ws.Range("D10").Value = 'Tính giới hạn \(\lim\limits_{x \to \infty} \frac{\sqrt{x^2 + 1} + x}{x + 1}\)'
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = ' \frac{1}{2}'
$ws.Range("H10").Value = ' +\infty'
$ws.Range("I10").Value = 2

$ws.Range("D11").Value = 'Tính giới hạn của dãy số sau khi \(n \to \infty\): \(x_n = \frac{1}{2} \left( x_{n-1} + \frac{1}{x_{n-1}} \right), \quad x_0 > 0\)'
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = ' +\infty'
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = ' \frac{1}{2}'
$ws.Range("I11").Value = 1

$ws.Range("D12").Value = 'Tính \(\lim_{x \to 2} \frac{x^2 - 4}{x - 2}\)'
$ws.Range("E12").Value = '\infty'
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 4

$ws.Range("D13").Value = 'Những VCB nào sau đây là tương đương:'
$ws.Range("E13").Value = '\arctan x , e^x - 1 , \frac{a^x - 1}{\ln a} , \ln(1 + x)'
$ws.Range("F13").Value = '\arcsin x , e^x , \frac{x^a - 1}{\ln x} , \ln(1 + x)'
$ws.Range("G13").Value = '\sin x , e^x , \frac{x^a - 1}{\ln x} , \ln(x)'
$ws.Range("H13").Value = '\cos x , e^x - 1 , \frac{a^x - 1}{\ln x} , \ln(1 + x)'
$ws.Range("I13").Value = '\arctan x , e^x - 1 , \frac{a^x - 1}{\ln a} , \ln(1 + x)'

# --- Step 5: the author's last actions left the view scrolled to / focused on I6 ---
$ws.Range("I6").Select()
$excel.ActiveWindow.ScrollRow = 3
